# Filmabrechnungen welche aus anderen Abrechnungsperioden entstanden sind
# werden von "Filmausgaben" nach "Sonstige Ausgaben" verschoben und dort
# als "Schulden aus anderen Abrechnungsperioden" kategorisiert.

$wb = $excel.ActiveWorkbook

$wsFilm = $wb.Worksheets.Item("Filmausgaben")
$wsSonst = $wb.Worksheets.Item("Sonstige Ausgaben")

# --- 1. Copy the two "Hotel Sinistra" rows (rows 2:3) from Filmausgaben
#        into the next free rows (9:10) of "Sonstige Ausgaben".
#        Sonstige Ausgaben has an extra "Kategorie" column (B), so the
#        remaining columns shift one to the right (A->A, B..H -> C..I).
$wsFilm.Range("A2:A3").Copy($wsSonst.Range("A9"))
$wsFilm.Range("B2:H3").Copy($wsSonst.Range("C9"))

# --- 2. Set the new "Kategorie" value for the moved rows.
$wsSonst.Range("B9").Value = "Schulden aus anderen Abrechnungsperioden "
$wsSonst.Range("B10").Value = "Schulden aus anderen Abrechnungsperioden "

# --- 3. Grow Table15 (Sonstige Ausgaben) so it covers the new rows.
$loSonst = $wsSonst.ListObjects.Item("Table15")
$loSonst.Resize($wsSonst.Range("A1:I10"))

# Column B ("Kategorie") now holds a much longer text, widen it and drop
# the old bestFit autosize.
$wsSonst.Columns.Item(2).ColumnWidth = 34.7

# --- 4. Remove the two rows from Filmausgaben (Table1 shrinks automatically).
$wsFilm.Range("A2:A3").EntireRow.Delete()

# --- 5. Restore a sensible selection on both affected sheets.
#        Filmausgaben stays the active (selected) tab, so activate it last.
$wsSonst.Activate()
[void]$wsSonst.Range("B13").Select()
$wsFilm.Activate()
[void]$wsFilm.Range("B13").Select()
